$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Numeric value fills (previously-empty Unit Price cells) ---
$ws.Range("C6").Value = 500
$ws.Range("C7").Value = 150
$ws.Range("C8").Value = 25
$ws.Range("C13").Value = 2
$ws.Range("C14").Value = 5
$ws.Range("C14").NumberFormat = '_($* #,##0.00_);_($* (#,##0.00);_($* "-"??_);_(@_)'
$ws.Range("C15").Value = 7
$ws.Range("C15").NumberFormat = '_($* #,##0.00_);_($* (#,##0.00);_($* "-"??_);_(@_)'
$ws.Range("C19").Value = 14.95

# --- 2. Text changes ---
$ws.Range("B15").Value = "40 inch strap with buckle"
$ws.Range("F14").Value = "approx"
$ws.Range("B19").Value = "FTDI Basic to load program onto the microcontroller"

# --- 3. Drop the old batchpcb link out of G9 (row 9 / Printed Wiring Board) ---
# (hyperlinks are handled as a full collection further below)
$ws.Range("G9").ClearContents()

# --- 4. Rebuild the hyperlinks in the target order (engine quirk: removing a
#     single Hyperlink object doesn't take effect, but it's easy to wipe the
#     whole collection and re-add every link in the right order) ---
$ws.Cells.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("G3"), "https://www.sparkfun.com/products/11114")
$ws.Hyperlinks.Add($ws.Range("G4"), "https://www.sparkfun.com/products/544?")
$ws.Hyperlinks.Add($ws.Range("G5"), "https://www.sparkfun.com/products/10160")
$ws.Hyperlinks.Add($ws.Range("G6"), "http://www.licor.com/env/products/light/quantum_sensors/190specs.html")
$ws.Hyperlinks.Add($ws.Range("G7"), "http://www.emesystems.com/uta_dat.htm")
$ws.Hyperlinks.Add($ws.Range("G10"), "https://www.sparkfun.com/products/8084")
$ws.Hyperlinks.Add($ws.Range("G12"), "https://www.sparkfun.com/products/117")
$ws.Hyperlinks.Add($ws.Range("G11"), "https://www.sparkfun.com/products/8235")
$ws.Hyperlinks.Add($ws.Range("G8"), "http://www.batteriesplus.com/product/32653-WKA6--8-dot2F-Werker-6V-8-dot2Ah-Battery/100085-1/102629-SLA-Sealed-Lead-Acid-Batteries/102647-Werker/6V.aspx")
$ws.Hyperlinks.Add($ws.Range("G20"), "http://www.amazon.com/SanDisk-microSD-Memory-Card-Adapter/dp/B0047WZOOO/ref=pd_bxgy_pc_text_z")
$ws.Hyperlinks.Add($ws.Range("G23"), "https://www.sparkfun.com/products/11367")
$ws.Hyperlinks.Add($ws.Range("G16"), "http://www.shopworldkitchen.com/snapware/rectangle-airtight-foodstorage-container-1098431")
$ws.Hyperlinks.Add($ws.Range("G15"), "http://www.rei.com/product/848429/redpoint-34-webbing-straps-with-side-release-buckles-package-of-2")

# --- 5. Restore original cell formatting on the hyperlink cells (Hyperlinks.Add
#     stamps its own style, so re-apply the workbook's own hyperlink styles -
#     "G13" is a plain never-linked s=7 cell, "G7" carries the left/vcenter s=6
#     variant used by the first two data rows) ---
$ws.Range("G13").Copy() | Out-Null
$ws.Range("G4,G5,G6,G8,G10,G11,G12,G15,G16,G20,G23").PasteSpecial(-4122) | Out-Null

$ws.Range("G7").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# --- 6. Selection moves to B25 ---
$ws.Range("B25").Select() | Out-Null
